$d = $word.ActiveDocument

# Locate the target paragraph: it currently holds the single-paragraph text
# "一项基于CAD平台的UML图绘制软件设计" together with the pPr (rFonts hint
# eastAsia) and the _GoBack bookmark.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*一项基于*CAD*平台的*UML*图绘制软件设计*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs.Item($d.Paragraphs.Count)
}

# Range covering only the paragraph's text (runs), excluding the trailing
# paragraph mark and the zero-width bookmark sitting right before it, so
# that InsertXML only replaces the run content and leaves the paragraph's
# own mark (pPr) + bookmark where they are.
$r = $target.Range
$r.MoveEnd(1, -1) | Out-Null

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>一项基于</w:t></w:r>
            <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>CAD</w:t></w:r>
            <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>平台的</w:t></w:r>
            <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>UML</w:t></w:r>
            <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>图绘制软件设计</w:t></w:r>
          </w:p>
          <w:p>
            <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>可以提供基本的</w:t></w:r>
            <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>UML</w:t></w:r>
            <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>图形设计，提供便捷的</w:t></w:r>
            <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>UML</w:t></w:r>
            <w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>图的输入与生成。同时支持代码生成与反向代码解析相关功能。</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>
            </w:pPr>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($xml) | Out-Null
